# Apply updated coin price/volume/name/link data (symbol-list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.80%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.54%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.819"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.37%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06272"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.81%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.866"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.81%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8770"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.85%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.9546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.96%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1456"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'4.05%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.05164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.25%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07279"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.75%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03131"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09049"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001548"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.18%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006279"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.62%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006092"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.48%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.459"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.25%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.269"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.00%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'5.02%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.61%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.00%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.860"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-5.62%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04318"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.74%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001176"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.03%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004276"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.64%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.17%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'2.74%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04042"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.51%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006701"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'62.22%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1153"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.62%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002131"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.63%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01411"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.63%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005162"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.81%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.19%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.304"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'3,342.74%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-12.18%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.19%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("E50").Style = "Normal"
